$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell content (values, formulas, shared strings) first -
# the sheet is being rebuilt with a different (smaller) data set and a new
# "sex" column.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "ref"
$ws.Range("B1").Value = "first_name__"
$ws.Range("C1").Value = "last_namess"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "entrance_datetime"
$ws.Range("F1").Value = "sex"

# Data row
$ws.Range("A2").Value = "test-"
$ws.Range("B2").Value = "c1c1c1c"
$ws.Range("C2").Value = "lol"
$ws.Range("D2").Value = "c2@gmail.com"
$ws.Range("E2").Value = "2023-01-01"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"
$ws.Range("F2").Value = "M"

# Column sizing: the wide column moves from D (firstName/email removed) to
# the new E (entrance_datetime) column.
$ws.Columns(5).ColumnWidth = 18.1
